$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 45741.01041666666
$ws.Cells.Item(2, 2).Value = 956
$ws.Cells.Item(3, 1).Value = 45741.02083333334
$ws.Cells.Item(3, 2).Value = 959
$ws.Cells.Item(4, 1).Value = 45741.03125
$ws.Cells.Item(4, 2).Value = 919
$ws.Cells.Item(5, 1).Value = 45741.04166666666
$ws.Cells.Item(5, 2).Value = 907
$ws.Cells.Item(6, 1).Value = 45741.05208333334
$ws.Cells.Item(6, 2).Value = 931
$ws.Cells.Item(7, 1).Value = 45741.0625
$ws.Cells.Item(7, 2).Value = 1003
$ws.Cells.Item(8, 1).Value = 45741.07291666666
$ws.Cells.Item(8, 2).Value = 999
$ws.Cells.Item(9, 1).Value = 45741.08333333334
$ws.Cells.Item(9, 2).Value = 992
$ws.Cells.Item(10, 1).Value = 45741.09375
$ws.Cells.Item(10, 2).Value = 982
$ws.Cells.Item(11, 1).Value = 45741.10416666666
$ws.Cells.Item(11, 2).Value = 996
$ws.Cells.Item(12, 1).Value = 45741.11458333334
$ws.Cells.Item(12, 2).Value = 1064
$ws.Cells.Item(13, 1).Value = 45741.125
$ws.Cells.Item(13, 2).Value = 1072
$ws.Cells.Item(14, 1).Value = 45741.13541666666
$ws.Cells.Item(14, 2).Value = 1061
$ws.Cells.Item(15, 1).Value = 45741.14583333334
$ws.Cells.Item(15, 2).Value = 1050
$ws.Cells.Item(16, 1).Value = 45741.15625
$ws.Cells.Item(16, 2).Value = 1014
$ws.Cells.Item(17, 1).Value = 45741.16666666666
$ws.Cells.Item(17, 2).Value = 1045
$ws.Cells.Item(18, 1).Value = 45741.17708333334
$ws.Cells.Item(18, 2).Value = 1061
$ws.Cells.Item(19, 1).Value = 45741.1875
$ws.Cells.Item(19, 2).Value = 1044
$ws.Cells.Item(20, 1).Value = 45741.19791666666
$ws.Cells.Item(20, 2).Value = 990
$ws.Cells.Item(21, 1).Value = 45741.20833333334
$ws.Cells.Item(21, 2).Value = 910
$ws.Cells.Item(22, 1).Value = 45741.21875
$ws.Cells.Item(22, 2).Value = 875
$ws.Cells.Item(23, 1).Value = 45741.22916666666
$ws.Cells.Item(23, 2).Value = 872
$ws.Cells.Item(24, 1).Value = 45741.23958333334
$ws.Cells.Item(24, 2).Value = 822
$ws.Cells.Item(25, 1).Value = 45741.25
$ws.Cells.Item(25, 2).Value = 829
$ws.Cells.Item(26, 1).Value = 45741.26041666666
$ws.Cells.Item(26, 2).Value = 847
$ws.Cells.Item(27, 1).Value = 45741.27083333334
$ws.Cells.Item(27, 2).Value = 811
$ws.Cells.Item(28, 1).Value = 45741.28125
$ws.Cells.Item(28, 2).Value = 820
$ws.Cells.Item(29, 1).Value = 45741.29166666666
$ws.Cells.Item(29, 2).Value = 831
$ws.Cells.Item(30, 1).Value = 45741.30208333334
$ws.Cells.Item(30, 2).Value = 752
$ws.Cells.Item(31, 1).Value = 45741.3125
$ws.Cells.Item(31, 2).Value = 724
$ws.Cells.Item(32, 1).Value = 45741.32291666666
$ws.Cells.Item(32, 2).Value = 727
$ws.Cells.Item(33, 1).Value = 45741.33333333334
$ws.Cells.Item(33, 2).Value = 762
$ws.Cells.Item(34, 1).Value = 45741.34375
$ws.Cells.Item(34, 2).Value = 750
$ws.Cells.Item(35, 1).Value = 45741.35416666666
$ws.Cells.Item(35, 2).Value = 762
$ws.Cells.Item(36, 1).Value = 45741.36458333334
$ws.Cells.Item(36, 2).Value = 810
$ws.Cells.Item(37, 1).Value = 45741.375
$ws.Cells.Item(37, 2).Value = 846
$ws.Cells.Item(38, 1).Value = 45741.38541666666
$ws.Cells.Item(38, 2).Value = 830
$ws.Cells.Item(39, 1).Value = 45741.39583333334
$ws.Cells.Item(39, 2).Value = 807
$ws.Cells.Item(40, 1).Value = 45741.40625
$ws.Cells.Item(40, 2).Value = 781
$ws.Cells.Item(41, 1).Value = 45741.41666666666
$ws.Cells.Item(41, 2).Value = 748
$ws.Cells.Item(42, 1).Value = 45741.42708333334
$ws.Cells.Item(42, 2).Value = 772
$ws.Cells.Item(43, 1).Value = 45741.4375
$ws.Cells.Item(43, 2).Value = 843
$ws.Cells.Item(44, 1).Value = 45741.44791666666
$ws.Cells.Item(44, 2).Value = 901
$ws.Cells.Item(45, 1).Value = 45741.45833333334
$ws.Cells.Item(45, 2).Value = 949
$ws.Cells.Item(46, 1).Value = 45741.46875
$ws.Cells.Item(47, 1).Value = 45741.47916666666
$ws.Cells.Item(48, 1).Value = 45741.48958333334
$ws.Cells.Item(49, 1).Value = 45741.5
$ws.Cells.Item(50, 1).Value = 45741.51041666666
$ws.Cells.Item(51, 1).Value = 45741.52083333334
$ws.Cells.Item(52, 1).Value = 45741.53125
$ws.Cells.Item(53, 1).Value = 45741.54166666666
$ws.Cells.Item(54, 1).Value = 45741.55208333334
$ws.Cells.Item(55, 1).Value = 45741.5625
$ws.Cells.Item(56, 1).Value = 45741.57291666666
$ws.Cells.Item(57, 1).Value = 45741.58333333334
$ws.Cells.Item(58, 1).Value = 45741.59375
$ws.Cells.Item(59, 1).Value = 45741.60416666666
$ws.Cells.Item(60, 1).Value = 45741.61458333334
$ws.Cells.Item(61, 1).Value = 45741.625
$ws.Cells.Item(62, 1).Value = 45741.63541666666
$ws.Cells.Item(63, 1).Value = 45741.64583333334
$ws.Cells.Item(64, 1).Value = 45741.65625
$ws.Cells.Item(65, 1).Value = 45741.66666666666
$ws.Cells.Item(66, 1).Value = 45741.67708333334
$ws.Cells.Item(67, 1).Value = 45741.6875
$ws.Cells.Item(68, 1).Value = 45741.69791666666
$ws.Cells.Item(69, 1).Value = 45741.70833333334
$ws.Cells.Item(70, 1).Value = 45741.71875
$ws.Cells.Item(71, 1).Value = 45741.72916666666
$ws.Cells.Item(72, 1).Value = 45741.73958333334
$ws.Cells.Item(73, 1).Value = 45741.75
$ws.Cells.Item(74, 1).Value = 45741.76041666666
$ws.Cells.Item(75, 1).Value = 45741.77083333334
$ws.Cells.Item(76, 1).Value = 45741.78125
$ws.Cells.Item(77, 1).Value = 45741.79166666666
$ws.Cells.Item(78, 1).Value = 45741.80208333334
$ws.Cells.Item(79, 1).Value = 45741.8125
$ws.Cells.Item(80, 1).Value = 45741.82291666666
$ws.Cells.Item(81, 1).Value = 45741.83333333334
$ws.Cells.Item(82, 1).Value = 45741.84375
$ws.Cells.Item(83, 1).Value = 45741.85416666666
$ws.Cells.Item(84, 1).Value = 45741.86458333334
$ws.Cells.Item(85, 1).Value = 45741.875
$ws.Cells.Item(86, 1).Value = 45741.88541666666
$ws.Cells.Item(87, 1).Value = 45741.89583333334
$ws.Cells.Item(88, 1).Value = 45741.90625
$ws.Cells.Item(89, 1).Value = 45741.91666666666
$ws.Cells.Item(90, 1).Value = 45741.92708333334
$ws.Cells.Item(91, 1).Value = 45741.9375
$ws.Cells.Item(92, 1).Value = 45741.94791666666
$ws.Cells.Item(93, 1).Value = 45741.95833333334
$ws.Cells.Item(94, 1).Value = 45741.96875
$ws.Cells.Item(95, 1).Value = 45741.97916666666
$ws.Cells.Item(96, 1).Value = 45741.98958333334
$ws.Cells.Item(97, 1).Value = 45742
